$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

# Dimension / institutional row: 0.06 -> 0.24
$r1 = $table.Cell(3, 3).Range
$r1.Find.Execute("0.06", $true, $false, $false, $false, $false, $true, 1, $false, "0.24", 1)

# Species / Cod row: 0.52 -> 0.34
$r2 = $table.Cell(5, 3).Range
$r2.Find.Execute("0.52", $true, $false, $false, $false, $false, $true, 1, $false, "0.34", 1)

# Species / Hake row: 0.06 -> 0.36
$r3 = $table.Cell(6, 3).Range
$r3.Find.Execute("0.06", $true, $false, $false, $false, $false, $true, 1, $false, "0.36", 1)
